$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A19 with the refreshed timestamp (tiny float precision change on re-retrieval)
$ws.Cells.Item(19, 1).Value = 44332.77772781134

# Append new row 20 with the newly retrieved data
$ws.Cells.Item(20, 1).Value = 44333.78536134069
$ws.Cells.Item(20, 2).Value = 73812
$ws.Cells.Item(20, 3).Value = 62235
$ws.Cells.Item(20, 4).Value = 3221
$ws.Cells.Item(20, 5).Value = 2077
$ws.Cells.Item(20, 6).Value = 1469
$ws.Cells.Item(20, 7).Value = 19104
$ws.Cells.Item(20, 8).Value = 1307
$ws.Cells.Item(20, 9).Value = 853
$ws.Cells.Item(20, 10).Value = 198
